$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Existing sheets: tweak selection / column width on "Sheet2" (the movies
#    sheet) and make it no longer the active tab (the new sheet takes over).
# ---------------------------------------------------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Sheet2: selection moves to C2, column C gets wider.
[void]$sheet2.Range("C2").Select()
$sheet2.Columns.Item(3).ColumnWidth = 15.15

# Sheet1 and Sheet2 tab colors re-saved as opaque white (alpha channel set).
$sheet2.Tab.Color = 16777215
$sheet1.Tab.Color = 16777215

# ---------------------------------------------------------------------------
# 2. Add the new worksheet "Movies with Blank Cells" as the 3rd / last sheet.
#    Adding it after the current last sheet also makes it the active tab,
#    matching the workbook's new activeTab="2".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Movies with Blank Cells"
$newSheet.Tab.Color = 16777215

# Header row.
$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "tit"
$newSheet.Range("C1").Value = "grs"
$newSheet.Range("D1").Value = "opn"

# Row 2 - fully populated.
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "Marvel's The Avengers"
$newSheet.Range("C2").Value = 623357910
$newSheet.Range("D2").Value = 41033
$newSheet.Range("D2").NumberFormat = "M/D/YYYY"

# Row 3 - grs (C3) left blank.
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "The Dark Knight Rises"
$newSheet.Range("D3").Value = 41110
$newSheet.Range("D3").NumberFormat = "M/D/YYYY"

# Row 4 - opn (D4) left blank (but still formatted as a date cell).
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "The Hunger Games"
$newSheet.Range("C4").Value = 408010692
$newSheet.Range("D4").NumberFormat = "M/D/YYYY"

# Row 5 - grs (C5) and opn (D5) both left blank.
$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = "Skyfall"
$newSheet.Range("D5").NumberFormat = "M/D/YYYY"

# Column widths for the new sheet.
$newSheet.Columns.Item(1).ColumnWidth = 8.83
$newSheet.Columns.Item(2).ColumnWidth = 19.33

# Select A2 as the active cell on the new sheet.
[void]$newSheet.Range("A2").Select()
